$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold, centered, bordered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for the new columns I (I0) and J (IF), rows 2-33
$data = @(
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(5, 6),
    @(10, 10),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(5, 6),
    @(6, 6),
    @(6, 6),
    @(5, 5),
    @(8, 8),
    @(7, 7),
    @(8, 9),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(8, 9),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
